$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 10: date, hours, description
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A10").Value = 42728
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "zuende übersetzt von allem was zu dem zeitpunkt existiert hat und benötigt wird/werden wird"

# Update selection to match the post-edit state
$ws.Range("H11").Select()
